# Updated cryptos list on Fri Apr  7 04:45:03 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as literal text (protects values like "313.03"
# or "  +0.18%  " from being reinterpreted as numbers/dates by Excel).
function Set-TextValue($cell, $text) {
    $ws.Range($cell).Value = "'" + $text
}

# Row 2 - Bitcoin
Set-TextValue "D2" "28.135.45"

# Row 3 - Ethereum
Set-TextValue "D3" "1.882.76"

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.18%  "

# Row 5 - BNB
Set-TextValue "D5" "313.03"
Set-TextValue "E5" "  -0.48%  "

# Row 6 - USDC
Set-TextValue "D6" "1.002"
Set-TextValue "E6" "  +0.14%  "

# Row 7 - XRP
Set-TextValue "D7" "0.5044"
Set-TextValue "E7" "  +0.43%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.3835"

# Row 9 - Dogecoin
Set-TextValue "D9" "0.08562"
Set-TextValue "E9" "  -7.35%  "

# Row 10 - Polygon
Set-TextValue "E10" "  -1.26%  "

# Row 11 - OKB
Set-TextValue "D11" "41.77"
Set-TextValue "E11" "  -0.25%  "

# Row 12 - Polkadot
Set-TextValue "D12" "6.272"
Set-TextValue "E12" "  -1.96%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.879.44"
Set-TextValue "E13" "  -1.23%  "

# Row 14 - Solana
Set-TextValue "D14" "20.59"
Set-TextValue "E14" "  -1.12%  "

# Row 15 - Chainlink
Set-TextValue "D15" "7.220"
Set-TextValue "E15" "  -0.94%  "

# Row 16 - BinanceUSD
Set-TextValue "E16" "  +0.20%  "

# Row 17 - ShibaInu
Set-TextValue "D17" "0.00001102"
Set-TextValue "E17" "  -0.92%  "

# Row 18 - Litecoin
Set-TextValue "D18" "91.28"
Set-TextValue "E18" "  -1.26%  "

# Row 19 - TRON
Set-TextValue "D19" "0.06669"
Set-TextValue "E19" "  +0.19%  "

# Row 20 - Avalanche
Set-TextValue "E20" "  +1.42%  "

# Row 21 - Dai
Set-TextValue "D21" "1.002"
Set-TextValue "E21" "  +0.22%  "

# Row 22 - Uniswap
Set-TextValue "D22" "6.106"
Set-TextValue "E22" "  -1.70%  "

# Row 23 - WrappedBTC
Set-TextValue "D23" "28.170.33"
Set-TextValue "E23" "  -0.07%  "

# Row 24 - Cosmos
Set-TextValue "D24" "11.19"
Set-TextValue "E24" "  -2.39%  "

# Row 25 - Toncoin
Set-TextValue "D25" "2.266"
Set-TextValue "E25" "  -2.34%  "

# Row 26 - LidoDAOToken
Set-TextValue "D26" "2.591"
Set-TextValue "E26" "  +1.38%  "

# Row 27 - WrappedliquidstakedEther2.0
Set-TextValue "D27" "2.093.08"
Set-TextValue "E27" "  -1.55%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "20.73"
Set-TextValue "E28" "  -0.67%  "

# Row 29 - Monero
Set-TextValue "D29" "156.35"
Set-TextValue "E29" "  -1.32%  "

# Row 30 - BitcoinCash
Set-TextValue "D30" "126.47"
Set-TextValue "E30" "  -0.41%  "

# Row 31 - Stellar
Set-TextValue "D31" "0.1053"
Set-TextValue "E31" "  -0.68%  "

# Row 32 - ImmutableX
Set-TextValue "D32" "1.056"
Set-TextValue "E32" "  -2.32%  "

# Row 33 - Filecoin
Set-TextValue "D33" "5.645"
Set-TextValue "E33" "  +0.64%  "

# Row 34 - HuobiToken
Set-TextValue "D34" "3.605"
Set-TextValue "E34" "  -0.39%  "

# Row 35 - FraxShare
Set-TextValue "D35" "9.717"
Set-TextValue "E35" "  +1.57%  "

# Row 36 - VeChain
Set-TextValue "D36" "0.02456"
Set-TextValue "E36" "  +2.37%  "

# Row 37 - Hedera
Set-TextValue "D37" "0.06560"
Set-TextValue "E37" "  -0.54%  "

# Row 38 - was ARBITRUM, now Algorand
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D38" "0.2180"
Set-TextValue "E38" "  -1.17%  "

# Row 39 - was Algorand, now ARBITRUM
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D39" "1.228"
Set-TextValue "E39" "  +0.47%  "

# Row 40 - TheSandbox
Set-TextValue "D40" "0.6585"
Set-TextValue "E40" "  +1.69%  "

# Row 41 - TrustWalletToken
Set-TextValue "D41" "1.244"
Set-TextValue "E41" "  -7.00%  "

# Row 42 - Aptos
Set-TextValue "D42" "11.41"
Set-TextValue "E42" "  +0.04%  "

# Row 43 - InternetComputer(DFINITY)
Set-TextValue "D43" "4.925"
Set-TextValue "E43" "  -1.06%  "

# Row 44 - Decentraland
Set-TextValue "D44" "0.6223"
Set-TextValue "E44" "  +1.89%  "

# Row 45 - EnergySwap
Set-TextValue "D45" "13.07"
Set-TextValue "E45" "  -2.19%  "

# Row 46 - WEMIXTOKEN
Set-TextValue "D46" "1.302"
Set-TextValue "E46" "  -0.31%  "

# Row 47 - PancakeSwap
Set-TextValue "D47" "3.687"
Set-TextValue "E47" "  -0.20%  "

# Row 48 - NEARProtocol
Set-TextValue "D48" "2.024"
Set-TextValue "E48" "  +1.02%  "

# Row 49 - EOS
Set-TextValue "E49" "  +1.28%  "

# Row 50 - Quant
Set-TextValue "D50" "121.08"
Set-TextValue "E50" "  -0.93%  "

# Row 51 - Aave
Set-TextValue "D51" "80.91"
Set-TextValue "E51" "  +2.32%  "
